# Auto-generated Excel COM-interop script
# Applies scheduled market-data refresh updates to the Leve profit sheets
# (columns H-N: currentAveragePrice*, LevePrice*, LeveProfit*)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 167.5
$ws.Range("I33").Value = 167.5
$ws.Range("K33").Value = 167.5
$ws.Range("M33").Value = 61.5

$ws.Range("H40").Value = 1413.909
$ws.Range("I40").Value = 1004.4286
$ws.Range("K40").Value = 1004.4286
$ws.Range("M40").Value = -829.4286

$ws.Range("H132").Value = 1744.4445
$ws.Range("I132").Value = 1818
$ws.Range("J132").Value = 825
$ws.Range("K132").Value = 5454
$ws.Range("L132").Value = 2475
$ws.Range("M132").Value = -2924
$ws.Range("N132").Value = -7535

$ws.Range("H137").Value = 1343.3513
$ws.Range("I137").Value = 1274.2963
$ws.Range("J137").Value = 1529.8
$ws.Range("K137").Value = 3822.8889
$ws.Range("L137").Value = 4589.4
$ws.Range("M137").Value = -1272.8889
$ws.Range("N137").Value = -9689.4

$ws.Range("H138").Value = 2475.2942
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 2475.2942
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 7425.882599999999
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -17705.8826

$ws.Range("H141").Value = 1394.579
$ws.Range("J141").Value = 3069.1428
$ws.Range("L141").Value = 9207.428400000001
$ws.Range("N141").Value = -19567.4284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1610.4865
$ws.Range("I61").Value = 1443.6
$ws.Range("J61").Value = 2325.7144
$ws.Range("K61").Value = 1443.6
$ws.Range("L61").Value = 2325.7144
$ws.Range("M61").Value = -1231.6
$ws.Range("N61").Value = -2749.7144

$ws.Range("H63").Value = 1500
$ws.Range("I63").Value = 1000
$ws.Range("K63").Value = 1000
$ws.Range("M63").Value = -314

$ws.Range("H66").Value = 1500
$ws.Range("I66").Value = 1000
$ws.Range("K66").Value = 5000
$ws.Range("M66").Value = -1568

$ws.Range("H74").Value = 33334050
$ws.Range("I74").Value = 47619544
$ws.Range("J74").Value = 1233.2222
$ws.Range("K74").Value = 47619544
$ws.Range("L74").Value = 1233.2222
$ws.Range("M74").Value = -47618670
$ws.Range("N74").Value = -2981.2222

$ws.Range("H77").Value = 33334050
$ws.Range("I77").Value = 47619544
$ws.Range("J77").Value = 1233.2222
$ws.Range("K77").Value = 238097720
$ws.Range("L77").Value = 6166.111
$ws.Range("M77").Value = -238093352
$ws.Range("N77").Value = -14902.111

$ws.Range("H132").Value = 13479.523
$ws.Range("I132").Value = 1424.7894
$ws.Range("K132").Value = 4274.3682
$ws.Range("M132").Value = -1744.3682

$ws.Range("H136").Value = 1610.4865
$ws.Range("I136").Value = 1443.6
$ws.Range("J136").Value = 2325.7144
$ws.Range("K136").Value = 4330.799999999999
$ws.Range("L136").Value = 6977.1432
$ws.Range("M136").Value = -1780.799999999999
$ws.Range("N136").Value = -12077.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1766.125
$ws.Range("I86").Value = 1446.2307
$ws.Range("J86").Value = 3152.3333
$ws.Range("K86").Value = 1446.2307
$ws.Range("L86").Value = 3152.3333
$ws.Range("M86").Value = -323.2307000000001
$ws.Range("N86").Value = -5398.3333

$ws.Range("H89").Value = 1766.125
$ws.Range("I89").Value = 1446.2307
$ws.Range("J89").Value = 3152.3333
$ws.Range("K89").Value = 7231.1535
$ws.Range("L89").Value = 15761.6665
$ws.Range("M89").Value = -1615.1535
$ws.Range("N89").Value = -26993.6665

$ws.Range("H94").Value = 2012
$ws.Range("I94").Value = 901.75
$ws.Range("J94").Value = 4232.5
$ws.Range("K94").Value = 901.75
$ws.Range("L94").Value = 4232.5
$ws.Range("M94").Value = -450.75
$ws.Range("N94").Value = -5134.5

$ws.Range("H134").Value = 4099.5
$ws.Range("I134").Value = 4529.778
$ws.Range("J134").Value = 1776
$ws.Range("K134").Value = 13589.334
$ws.Range("L134").Value = 5328
$ws.Range("M134").Value = -11054.334
$ws.Range("N134").Value = -10398

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12237.543
$ws.Range("I31").Value = 25714.77
$ws.Range("J31").Value = 4273.727
$ws.Range("K31").Value = 25714.77
$ws.Range("L31").Value = 4273.727
$ws.Range("M31").Value = -25419.77
$ws.Range("N31").Value = -4863.727

$ws.Range("H34").Value = 12237.543
$ws.Range("I34").Value = 25714.77
$ws.Range("J34").Value = 4273.727
$ws.Range("K34").Value = 25714.77
$ws.Range("L34").Value = 4273.727
$ws.Range("M34").Value = -25512.77
$ws.Range("N34").Value = -4677.727

$ws.Range("H58").Value = 9492.648999999999
$ws.Range("I58").Value = 668.9
$ws.Range("K58").Value = 668.9
$ws.Range("M58").Value = -465.9

$ws.Range("H122").Value = 1194.7222
$ws.Range("I122").Value = 1180.909
$ws.Range("J122").Value = 1216.4286
$ws.Range("K122").Value = 3542.727
$ws.Range("L122").Value = 3649.2858
$ws.Range("M122").Value = -1092.727
$ws.Range("N122").Value = -8549.2858

$ws.Range("H132").Value = 16669.111
$ws.Range("I132").Value = 23043.709
$ws.Range("J132").Value = 3919.9167
$ws.Range("K132").Value = 69131.12699999999
$ws.Range("L132").Value = 11759.7501
$ws.Range("M132").Value = -66601.12699999999
$ws.Range("N132").Value = -16819.7501

$ws.Range("H134").Value = 734.913
$ws.Range("I134").Value = 649.95
$ws.Range("J134").Value = 1301.3334
$ws.Range("K134").Value = 1949.85
$ws.Range("L134").Value = 3904.0002
$ws.Range("M134").Value = 585.1499999999999
$ws.Range("N134").Value = -8974.0002

$ws.Range("H135").Value = 50600
$ws.Range("J135").Value = 50600
$ws.Range("L135").Value = 50600
$ws.Range("N135").Value = -60740

$ws.Range("H136").Value = 9492.648999999999
$ws.Range("I136").Value = 668.9
$ws.Range("K136").Value = 2006.7
$ws.Range("M136").Value = 543.3000000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 959.2857
$ws.Range("J5").Value = 801.6667
$ws.Range("L5").Value = 2405.0001
$ws.Range("N5").Value = -2629.0001

$ws.Range("H107").Value = 4111.5386
$ws.Range("I107").Value = 16881.666
$ws.Range("K107").Value = 50644.99800000001
$ws.Range("M107").Value = -48724.99800000001

$ws.Range("H122").Value = 455.25
$ws.Range("I122").Value = 234.88889
$ws.Range("K122").Value = 2114.00001
$ws.Range("M122").Value = 335.9999899999998

$ws.Range("H131").Value = 760.1
$ws.Range("J131").Value = 775.3711499999999
$ws.Range("L131").Value = 2326.11345
$ws.Range("N131").Value = -12406.11345

$ws.Range("H132").Value = 1163.9048
$ws.Range("I132").Value = 499.5
$ws.Range("K132").Value = 4495.5
$ws.Range("M132").Value = -1965.5

$ws.Range("H135").Value = 959.2857
$ws.Range("J135").Value = 801.6667
$ws.Range("L135").Value = 7215.0003
$ws.Range("N135").Value = -12285.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 20613.066
$ws.Range("I132").Value = 4336.091
$ws.Range("J132").Value = 65374.75
$ws.Range("K132").Value = 13008.273
$ws.Range("L132").Value = 196124.25
$ws.Range("M132").Value = -10478.273
$ws.Range("N132").Value = -201184.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 436.9375
$ws.Range("I16").Value = 436.9375
$ws.Range("K16").Value = 436.9375
$ws.Range("M16").Value = -266.9375

$ws.Range("H88").Value = 38000
$ws.Range("J88").Value = 38000
$ws.Range("L88").Value = 38000
$ws.Range("N88").Value = -38856

$ws.Range("H91").Value = 38000
$ws.Range("J91").Value = 38000
$ws.Range("L91").Value = 38000
$ws.Range("N91").Value = -40964

$ws.Range("H132").Value = 1871.4517
$ws.Range("I132").Value = 1330.1111
$ws.Range("J132").Value = 2621
$ws.Range("K132").Value = 3990.3333
$ws.Range("L132").Value = 7863
$ws.Range("M132").Value = -1460.3333
$ws.Range("N132").Value = -12923

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 5051326.5
$ws.Range("I107").Value = 1032.5
$ws.Range("J107").Value = 15151915
$ws.Range("K107").Value = 3097.5
$ws.Range("L107").Value = 45455745
$ws.Range("M107").Value = -1177.5
$ws.Range("N107").Value = -45459585

$ws.Range("H132").Value = 877.625
$ws.Range("I132").Value = 609.97437
$ws.Range("K132").Value = 1829.92311
$ws.Range("M132").Value = 700.0768899999998
